$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Record attendance (PM, Eng#1, Eng#2 present) for the Dec 3, 2018 meeting,
# which is row 29 of the attendance table.
$ws.Range("C29:E29").Value = 1

# Scroll the sheet down and move the active selection, matching the saved
# view state of the workbook (topLeftCell A7 -> A16, selection C29 -> H29).
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("H29").Select()
